# SimPathsEU parameters.xlsx edit
# 1) Add SAVINGS_RATE parameter row to the "Parameters" sheet.
# 2) Re-style the numeric tax/benefit-regime block (rows 23-36) to a plain font.
# 3) Add a new "Info" sheet describing every parameter, with the SAVINGS_RATE
#    description appended at the end.

$wb = $excel.ActiveWorkbook
$params = $wb.Worksheets.Item(1)
$params.Name = "Parameters"

# ---------------------------------------------------------------------------
# 1) New SAVINGS_RATE row at the bottom of Parameters (row 37)
# ---------------------------------------------------------------------------
$params.Range("A37").Value = "SAVINGS_RATE"
$params.Range("B37").Value = 0.056

$params.Range("A37:B37").Font.Name = "Aptos Narrow"
$params.Range("A37:B37").Font.Size = 12
$params.Range("A37").RowHeight = 16

# ---------------------------------------------------------------------------
# 2) Re-style rows 23-36 (key + value) to the plain Helvetica Neue font
# ---------------------------------------------------------------------------
$params.Range("A23:B36").Font.Name = "Helvetica Neue"
$params.Range("A23:B36").Font.Size = 10
$params.Range("A23:B36").Font.Bold = $false

# Selection / view bookkeeping on the Parameters sheet
$params.Range("A1:B37").Select()

# ---------------------------------------------------------------------------
# 3) New "Info" sheet
# ---------------------------------------------------------------------------
$info = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $params)
$info.Name = "Info"

$info.Columns.Item(1).ColumnWidth = 57.5
$info.Columns.Item(2).ColumnWidth = 93.33203125

$info.Range("A1").Value = "This file is used to rewrite the following country-specific parameters"
$info.Range("A1").Font.Name = "Aptos Narrow"
$info.Range("A1").Font.Size = 12
$info.Range("A1").RowHeight = 16

$info.Range("A2").Font.Name = "Aptos Narrow"
$info.Range("A2").Font.Size = 12
$info.Range("B2").Font.Name = "Aptos Narrow"
$info.Range("B2").Font.Size = 12
$info.Range("A2").RowHeight = 16

$info.Range("A3").Value = "Parameter Name"
$info.Range("B3").Value = "Description"
$info.Range("A3:B3").Font.Name = "Aptos Narrow"
$info.Range("A3:B3").Font.Size = 12
$info.Range("A3:B3").Font.Bold = $true
$info.Range("A3").RowHeight = 16

$keys = @(
    "MIN_AGE_TO_HAVE_INCOME",
    "MAX_LABOUR_HOURS_IN_WEEK",
    "HOURS_IN_WEEK",
    "USE_CONTINUOUS_LABOUR_SUPPLY_HOURS",
    "AGE_TO_BECOME_RESPONSIBLE",
    "MIN_AGE_TO_LEAVE_EDUCATION",
    "MAX_AGE_TO_LEAVE_CONTINUOUS_EDUCATION",
    "MAX_AGE_TO_ENTER_EDUCATION",
    "MIN_AGE_TO_RETIRE",
    "DEFAULT_AGE_TO_RETIRE",
    "MIN_AGE_FORMAL_SOCARE",
    "MIN_AGE_FLEXIBLE_LABOUR_SUPPLY",
    "MAX_AGE_FLEXIBLE_LABOUR_SUPPLY",
    "SHARE_OF_WEALTH_TO_ANNUITISE_AT_RETIREMENT",
    "ANNUITY_RATE_OF_RETURN",
    "MIN_HOURS_FULL_TIME_EMPLOYED",
    "MIN_HOURLY_WAGE_RATE",
    "MAX_HOURLY_WAGE_RATE",
    "MAX_HOURS_WEEKLY_FORMAL_CARE",
    "MAX_HOURS_WEEKLY_INFORMAL_CARE",
    "CHILDCARE_COST_EARNINGS_CAP",
    "TAXDB_REGIMES",
    "MIN_START_YEAR",
    "MAX_START_YEAR",
    "MIN_START_YEAR_TRAINING",
    "MAX_START_YEAR_TRAINING",
    "MIN_CAPITAL_INCOME_PER_MONTH",
    "MAX_CAPITAL_INCOME_PER_MONTH",
    "MIN_PERSONAL_PENSION_PER_MONTH",
    "MAX_PERSONAL_PENSION_PER_MONTH",
    "MAX_CHILD_AGE_FOR_FORMAL_CARE",
    "MIN_AGE_MATERNITY",
    "MAX_AGE_MATERNITY",
    "BASE_PRICE_YEAR",
    "PROB_NEWBORN_IS_MALE"
)

$descriptions = @(
    "Minimum age to have non-employment, non-benefit income",
    "Maximum number of labour hours allowed in a week",
    "Total hours in a week (used to calculate leisure in labour supply)",
    "If true, generates random weekly labour supply hours within each bracket; if false, uses fixed hours for all persons",
    "Age at which a person becomes reference person of their own benefit unit",
    "Minimum age to leave full-time education",
    "Maximum age to remain in continuous education",
    "Maximum age to enter education",
    "Minimum age to consider retirement",
    "Default retirement age (if pension included but retirement decision not modeled)",
    "Minimum age to receive formal social care",
    "Minimum age for flexible labour supply eligibility",
    "Maximum age for flexible labour supply eligibility",
    "Proportion of wealth to annuitise at retirement",
    "Assumed annuity rate of return",
    "Minimum weekly hours defining full-time employment",
    "Minimum possible hourly wage",
    "Maximum possible hourly wage",
    "Maximum number of hours of formal care per week",
    "Maximum number of hours of informal care per week",
    "Maximum share of earnings payable as childcare",
    "Number of tax/benefit regimes supported",
    "Minimum allowed simulation start year (oldest initial population)",
    "Maximum allowed simulation start year (most recent initial population)",
    "Minimum allowed training start year",
    "Maximum allowed training start year",
    "Minimum capital income per month",
    "Maximum capital income per month",
    "Minimum pension income per month",
    "Maximum pension income per month",
    "Maximum age of child eligible for formal care",
    "Minimum age a person can give birth",
    "Maximum age a person can give birth",
    "Base year for model parameters (prices)",
    "Probability a newborn is male"
)

$row = 4
for ($i = 0; $i -lt $keys.Length; $i++) {
    $info.Range("A$row").Value = $keys[$i]
    $info.Range("B$row").Value = $descriptions[$i]
    $info.Range("A$row").Font.Name = "Arial Unicode MS"
    $info.Range("A$row").Font.Size = 10
    $info.Range("B$row").Font.Name = "Aptos Narrow"
    $info.Range("B$row").Font.Size = 12
    $info.Range("A$row").RowHeight = 17
    $row++
}

$info.Range("A39").Value = "SAVINGS_RATE"
$info.Range("B39").Value = "Country-specific savings rate"
$info.Range("A39:B39").Font.Name = "Aptos Narrow"
$info.Range("A39:B39").Font.Size = 12
$info.Range("A39").RowHeight = 16

$info.Range("A1:B39").Select()
